# Update countries & provincias Spain
# Refresh case numbers for a handful of countries and update the
# "last updated" timestamp. Two pairs of adjacent rows (Kuwait /
# Emiratos Arabes Unidos and Estonia / Islandia) swap order because the
# refreshed totals change their rank when the sheet is kept sorted by
# "Casos totales" (column B) descending, so those rows' country-name
# cells are updated as well as their numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Last updated timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 14:24"

# --- India (row 5) ---
$ws.Range("B5").Value = 7552262
$ws.Range("C5").Value = 4024
$ws.Range("E5").Value = 773984
$ws.Range("G5").Value = 28
$ws.Range("H5").Value = 114670

# --- Irak (row 18) ---
$ws.Range("B18").Value = 430678
$ws.Range("C18").Value = 4044
$ws.Range("D18").Value = 363532
$ws.Range("E18").Value = 56829
$ws.Range("G18").Value = 63
$ws.Range("H18").Value = 10317

# --- Row 42: now Emiratos Arabes Unidos (was Kuwait) ---
$ws.Range("A42").Value = "Emiratos Arabes Unidos"
$ws.Range("B42").Value = 116517
$ws.Range("C42").Value = 915
$ws.Range("D42").Value = 108811
$ws.Range("E42").Value = 7240
$ws.Range("G42").Value = 3
$ws.Range("H42").Value = 466

# --- Row 43: now Kuwait (was Emiratos Arabes Unidos) ---
$ws.Range("A43").Value = "Kuwait"
$ws.Range("B43").Value = 116146
$ws.Range("D43").Value = 107860
$ws.Range("E43").Value = 7585
$ws.Range("H43").Value = 701

# --- Suiza (row 57) ---
$ws.Range("D57").Value = 53400
$ws.Range("E57").Value = 27627

# --- Libia (row 70) ---
$ws.Range("B70").Value = 49949
$ws.Range("C70").Value = 1159
$ws.Range("D70").Value = 27262
$ws.Range("E70").Value = 21955
$ws.Range("G70").Value = 7
$ws.Range("H70").Value = 732

# --- Azerbaiyan (row 74) ---
$ws.Range("B74").Value = 45295
$ws.Range("C74").Value = 331
$ws.Range("D74").Value = 40144
$ws.Range("E74").Value = 4521
$ws.Range("G74").Value = 4
$ws.Range("H74").Value = 630

# --- Dinamarca (row 81) ---
$ws.Range("B81").Value = 35844
$ws.Range("C81").Value = 452
$ws.Range("D81").Value = 29562
$ws.Range("E81").Value = 5596
$ws.Range("G81").Value = 6
$ws.Range("H81").Value = 686

# --- Bosnia y Herzegovina (row 82) ---
$ws.Range("B82").Value = 34661
$ws.Range("C82").Value = 549
$ws.Range("D82").Value = 25442
$ws.Range("E82").Value = 8222
$ws.Range("G82").Value = 13
$ws.Range("H82").Value = 997

# --- Row 141: now Islandia (was Estonia) ---
$ws.Range("A141").Value = "Islandia"
$ws.Range("B141").Value = 4101
$ws.Range("C141").Value = 46
$ws.Range("D141").Value = 2856
$ws.Range("E141").Value = 1234
$ws.Range("H141").Value = 11

# --- Row 142: now Estonia (was Islandia) ---
$ws.Range("A142").Value = "Estonia"
$ws.Range("B142").Value = 4085
$ws.Range("C142").Value = 7
$ws.Range("D142").Value = 3229
$ws.Range("E142").Value = 788
$ws.Range("H142").Value = 68

# --- Somalia (row 144) ---
$ws.Range("B144").Value = 3890
$ws.Range("C144").Value = 26
$ws.Range("E144").Value = 702

# --- Vietnam (row 168) ---
$ws.Range("B168").Value = 1140
$ws.Range("C168").Value = 6
$ws.Range("D168").Value = 1046
$ws.Range("E168").Value = 59

# --- Barbados (row 192) ---
$ws.Range("B192").Value = 222
$ws.Range("C192").Value = 1
$ws.Range("D192").Value = 203
$ws.Range("E192").Value = 12
